# Replace curly double quotes (“ ”) with straight single quotes (')
# in the English (en_US, column C) dialogue lines of the story sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$leftDoubleQuote  = [char]0x201C   # “
$rightDoubleQuote = [char]0x201D   # ”
$straightQuote    = "'"

# Cells in column C (en_US) that contain the curly double quotes to fix.
$targetCells = @("C3", "C6", "C18", "C21", "C23", "C46", "C73", "C78", "C97", "C101")

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    $text = $text.Replace($leftDoubleQuote, $straightQuote)
    $text = $text.Replace($rightDoubleQuote, $straightQuote)
    $cell.Value = $text
}
